# feat: add 2022-Q3 data
#
# - Insert a new "2022-Q3" sheet right after "总计", pushing 2022-Q2 / 2022-Q1 /
#   2021-Q4 one position later (tab order becomes:
#   总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4).
# - Prepend a 2022-Q3 row to the "总计" summary sheet and renumber the index
#   column so it again reads 0,1,2,3 top to bottom.
# - Populate the new "2022-Q3" sheet with the per-fund holding breakdown,
#   matching the layout used by the other quarter sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q3 and
#    shift the previous rows down.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("A2").Font.Bold = $true
$wsTotal.Range("A2").HorizontalAlignment = -4108
$wsTotal.Range("A2").VerticalAlignment = -4160
$wsTotal.Range("A2").Borders.LineStyle = 1

$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 1.42
$wsTotal.Range("B2:D2").Style = "Normal"

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add($null, $wsTotal)
$wsQ3.Name = "2022-Q3"

# Header row.
$header = $wsQ3.Range("B1:H1")
$header.Value = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

function Add-FundRow($row, $idx, $code, $name, $scale, $position, $ratio, $value, $rank) {
    $wsQ3.Range("A$row").Value = $idx
    $wsQ3.Range("A$row").Font.Bold = $true
    $wsQ3.Range("A$row").HorizontalAlignment = -4108
    $wsQ3.Range("A$row").VerticalAlignment = -4160
    $wsQ3.Range("A$row").Borders.LineStyle = 1

    $textRange = $wsQ3.Range("B$row`:G$row")
    $textRange.NumberFormat = "@"
    $wsQ3.Range("B$row").Value = $code
    $wsQ3.Range("C$row").Value = $name
    $wsQ3.Range("D$row").Value = $scale
    $wsQ3.Range("E$row").Value = $position
    $wsQ3.Range("F$row").Value = $ratio
    $wsQ3.Range("G$row").Value = $value
    $textRange.Style = "Normal"

    $wsQ3.Range("H$row").Value = $rank
}

Add-FundRow 2 0 "398021" "中海能源策略混合"             "21.73" "88.03" "4.37" "0.9496" 7
Add-FundRow 3 1 "160926" "大成创业板两年定期开放混合A"   "7.68"  "64.14" "2.57" "0.1974" 10
Add-FundRow 4 2 "398061" "中海消费混合"                 "3.91"  "85.30" "4.21" "0.1646" 8
Add-FundRow 5 3 "009798" "大成创业板两年定期开放混合C"   "2.71"  "64.14" "2.57" "0.0696" 10
Add-FundRow 6 4 "000166" "中海信息产业精选混合"          "0.77"  "89.31" "4.84" "0.0373" 5
Add-FundRow 7 5 "970073" "东证融汇成长优选混合A"         "0.37"  "91.31" "0.76" "0.0028" 3
Add-FundRow 8 6 "970074" "东证融汇成长优选混合C"         "0.12"  "91.31" "0.76" "0.0009" 3

# ---------------------------------------------------------------------------
# 3. Restore "总计" as the active/displayed sheet (adding a sheet above made
#    the new one active).
# ---------------------------------------------------------------------------
$wsTotal.Activate()
